# Add other computations to arianna spreadsheet
#
# - Fix a typo in the D1 header: "Vomega_eff(cm^2 sr)" -> "Vomega_eff(cm^3 sr)"
# - Append a new computation block (rows 20-30): headers in row 20 (B:D) and
#   a 10-row table (rows 21-30) reusing the Lint(cm) values already computed
#   in column F of the first block (rows 3-12), converting to lin(km) and
#   multiplying by 10^A to get a new aomega_eff(km^2 sr)-style quantity.
# - Leave the final selection on C31 (just past the new table), matching
#   where the author's cursor ended up after typing the data in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix header typo -------------------------------------------------
$ws.Range("D1").Value = "Vomega_eff(cm^3 sr)"

# --- new block header (row 20) ---------------------------------------
$ws.Range("B20").Value = "aomega_eff(km^2 sr)"
$ws.Range("C20").Value = "lint(cm)"
$ws.Range("D20").Value = "lin(km)"

# --- new data rows (21-30) --------------------------------------------
$rows = @(
  @{ r=21; a=-6.424895917;  c=161293479.880792 },
  @{ r=22; a=-4.848875917;  c=106198115.902776 },
  @{ r=23; a=-3.819635917;  c=69922478.1413035 },
  @{ r=24; a=-3.047705917;  c=46038038.5081131 },
  @{ r=25; a=-2.404435917;  c=30312154.9180693 },
  @{ r=26; a=-1.857645917;  c=19957990.512891  },
  @{ r=27; a=-1.439515917;  c=13140648.9043511 },
  @{ r=28; a=-1.021395917;  c=8652005.99809334 },
  @{ r=29; a=-0.6997559174; c=5696614.24910735 },
  @{ r=30; a=-0.5067759174; c=3750738.72004761 }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Range("A$r").Value = $row.a
  $ws.Range("B$r").Formula = "=10^A$r"
  $ws.Range("C$r").Value = $row.c
  $ws.Range("D$r").Formula = "=C$r/100000"
  $ws.Range("E$r").Formula = "=D$r*B$r"
}

# --- leave selection where the author's cursor ended up ---------------
$ws.Range("C31").Select()
